$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates as Excel serial numbers, matching existing column A values)
$data = @(
    @{ Row = 234; A = 44308; B = 3; C = 4; D = 249.3765586034913 },
    @{ Row = 235; A = 44309; B = 0; C = 4; D = 249.3765586034913 },
    @{ Row = 236; A = 44310; B = 0; C = 4; D = 249.3765586034913 },
    @{ Row = 237; A = 44311; B = 0; C = 3; D = 187.0324189526185 },
    @{ Row = 238; A = 44312; B = 0; C = 3; D = 187.0324189526185 }
)

foreach ($entry in $data) {
    $r = $entry.Row

    # Column A: date serial value, styled like the existing date column (copy full format from A233)
    $ws.Cells.Item($r - 1, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 1).Value = $entry.A

    # Columns B, C, D: plain numeric values
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
}
$excel.CutCopyMode = 0
